$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("libraries")

# --- libraries sheet -------------------------------------------------
# Drop iText (AGPL-licensed PDF export) and fold Batik into the new
# "Apache Batik" row; "Apache FOP" takes over PDF export duties. The
# remaining rows keep their existing order, so everything between the
# edited rows shifts up by one.

$ws.Range("A3").Value = 'Apache Batik'
$ws.Range("B3").Value = 'exporting images to SVG, EPS, …'
$ws.Range("C3").Value = 'https://xmlgraphics.apache.org/batik/'
$ws.Range("D3").Formula = "=""1.17"""
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = 'Apache 2.0'

$ws.Range("A4").Value = 'Apache FOP'
$ws.Range("B4").Value = 'exporting displayed graphs to PDF format'
$ws.Range("D4").Formula = "=""2.9"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("C4").Value = 'https://xmlgraphics.apache.org/fop/'
$ws.Range("E4").Value = 'Apache 2.0'

$ws.Range("A5").Value = 'Args4J'
$ws.Range("B5").Value = 'command-line option parsing'
$ws.Range("C5").Value = 'http://args4j.kohsuke.org/'
$ws.Range("D5").Value = '2.0.26'
$ws.Range("E5").Value = 'MIT'

$ws.Range("A6").Value = 'JAXB'
$ws.Range("B6").Value = 'access to an XML document from a Java program'
$ws.Range("C6").Value = 'https://javaee.github.io/jaxb-v2/'
$ws.Range("D6").Value = '3.0.2'
$ws.Range("E6").Value = 'CDDL 1.1 or GPL 2.0 with classpath exception'

$ws.Range("A7").Value = 'EMF'
$ws.Range("B7").Value = 'converting to and from Eclipse ecore format'
$ws.Range("C7").Value = 'http://eclipse.org'
$ws.Range("D7").Value = '2.7.0'
$ws.Range("E7").Value = 'EPL 2.0'

$ws.Range("A8").Value = 'Gnu Prolog'
$ws.Range("B8").Value = 'interpreting Prolog queries'
$ws.Range("C8").Value = 'http://www.gnu.org/software/gnuprologjava/'
$ws.Range("D8").Value = '0.2.6'
$ws.Range("E8").Value = 'LGPL 3.0'

$ws.Range("A9").Value = 'Groovy'
$ws.Range("B9").Value = 'easy and flexible access to the GROOVE API'
$ws.Range("C9").Value = 'http://groovy.codehaus.org/'
$ws.Range("D9").Value = '2.0.5'
$ws.Range("E9").Value = 'Apache 2.0'

$ws.Range("A10").Value = 'Jakarta activation'
$ws.Range("B10").Value = 'dynamic data manipulation'
$ws.Range("C10").Value = 'https://projects.eclipse.org/projects/ee4j.jaf'
$ws.Range("D10").Value = '3.0.1'
$ws.Range("E10").Value = 'EDL 1.0'

$ws.Range("A11").Value = 'Jakarta bind'
$ws.Range("B11").Value = 'mapping between XML documents and Java objects'
$ws.Range("C11").Value = 'https://jakarta.ee/specifications/xml-binding/'
$ws.Range("D11").Value = '1.2.2'
$ws.Range("E11").Value = 'EDL 1.0'

$ws.Range("A12").Value = 'Java Annotations'
$ws.Range("B12").Value = 'runtime analysis of Java annotations'
$ws.Range("C12").Value = 'https://wiki.eclipse.org/JDT_Core/Null_Analysis'
$ws.Range("D12").Value = '2.2.700'
$ws.Range("E12").Value = 'EPL 2.0'

$excel.CutCopyMode = 0

# Recalculate the workbook so the mirrored "save to libraries.csv" sheet
# (which references libraries!A2:E20 by formula) picks up the new rows.
$excel.CalculateFullRebuild()

# --- restore the end-user selections left behind by the edit ---------
$ws.Rows.Item(4).Select()

$ws2 = $wb.Worksheets.Item("save to libraries.csv")
$ws2.Activate()
$ws2.Range("C21").Select()
